# Remove slide 9 ("M13: Het project gebruikt ISO-25010 voor de specificatie
# van productkwaliteitseisen"). This measure has been deprecated; its
# content (use of ISO-25010) is already covered by M01.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$s.Delete()
